# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 542
    $ws.Range("F3").Value = 3501
    $ws.Range("F4").Value = 96
    $ws.Range("F5").Value = 683
}
